$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.32446187183548
$ws.Range("C2").Value = 8.112591183731475
$ws.Range("D2").Value = 12.80886201729186
$ws.Range("E2").Value = 12.49517531015564
$ws.Range("G2").Value = 62.74941152461878
$ws.Range("H2").Value = 22.66345704985719
$ws.Range("I2").Value = 35.0819008788039
$ws.Range("J2").Value = 7.82328645114011
$ws.Range("L2").Value = 13.25226449292265
$ws.Range("M2").Value = 19.8944104374259
$ws.Range("B3").Value = 21.05514656555827
$ws.Range("C3").Value = 7.732969483734005
$ws.Range("D3").Value = 12.82569295084351
$ws.Range("E3").Value = 12.51664045950034
$ws.Range("G3").Value = 62.66108474500914
$ws.Range("H3").Value = 22.6990681736813
$ws.Range("I3").Value = 35.15168699255292
$ws.Range("J3").Value = 7.809683811305415
$ws.Range("L3").Value = 13.26644612982795
$ws.Range("M3").Value = 19.854263752034
$ws.Range("B4").Value = 20.89325444276105
$ws.Range("C4").Value = 7.492153797742179
$ws.Range("D4").Value = 12.83800291574371
$ws.Range("E4").Value = 12.53057663289421
$ws.Range("G4").Value = 62.62349514967348
$ws.Range("H4").Value = 22.72522180382967
$ws.Range("I4").Value = 35.2014028153942
$ws.Range("J4").Value = 7.801111214766101
$ws.Range("L4").Value = 13.27705116190803
$ws.Range("M4").Value = 19.83322077822097
$ws.Range("B5").Value = 20.82823333698515
$ws.Range("C5").Value = 7.392227137434919
$ws.Range("D5").Value = 12.84351593939942
$ws.Range("E5").Value = 12.53644645946455
$ws.Range("G5").Value = 62.61236139045971
$ws.Range("H5").Value = 22.73695523306187
$ws.Range("I5").Value = 35.22338350226163
$ws.Range("J5").Value = 7.797562341402006
$ws.Range("L5").Value = 13.28185012622998
$ws.Range("M5").Value = 19.82555891883084
$ws.Range("B6").Value = 20.81749636086799
$ws.Range("C6").Value = 7.375530963924091
$ws.Range("D6").Value = 12.84446135978324
$ws.Range("E6").Value = 12.53743267355703
$ws.Range("G6").Value = 62.61076526088497
$ws.Range("H6").Value = 22.738968440947
$ws.Range("I6").Value = 35.22713715718968
$ws.Range("J6").Value = 7.796969702798175
$ws.Range("L6").Value = 13.28267582433171
$ws.Range("M6").Value = 19.82434200436262
$ws.Range("B7").Value = 20.89237359071476
$ws.Range("C7").Value = 7.490813194074971
$ws.Range("D7").Value = 12.83807525585408
$ws.Range("E7").Value = 12.53065502253439
$ws.Range("G7").Value = 62.62332805794729
$ws.Range("H7").Value = 22.72537569367446
$ws.Range("I7").Value = 35.20169229381116
$ws.Range("J7").Value = 7.801063578040206
$ws.Range("L7").Value = 13.27711394954098
$ws.Range("M7").Value = 19.83311374199132
$ws.Range("B8").Value = 21.23092996649845
$ws.Range("C8").Value = 7.9833850724055
$ws.Range("D8").Value = 12.81425517030226
$ws.Range("E8").Value = 12.50241978565299
$ws.Range("G8").Value = 62.71550275068529
$ws.Range("H8").Value = 22.67484417382032
$ws.Range("I8").Value = 35.10453460040974
$ws.Range("J8").Value = 7.818642016539829
$ws.Range("L8").Value = 13.25676064450275
$ws.Range("M8").Value = 19.879823126264
$ws.Range("B9").Value = 21.91881941562163
$ws.Range("C9").Value = 8.882426193638848
$ws.Range("D9").Value = 12.78322718713144
$ws.Range("E9").Value = 12.45303024554167
$ws.Range("G9").Value = 63.02818761756927
$ws.Range("H9").Value = 22.60989092536351
$ws.Range("I9").Value = 34.96873661515928
$ws.Range("J9").Value = 7.851365397165297
$ws.Range("L9").Value = 13.23189296051088
$ws.Range("M9").Value = 19.9997392738207
$ws.Range("B10").Value = 22.43404326371508
$ws.Range("C10").Value = 9.495971081509785
$ws.Range("D10").Value = 12.76999830429415
$ws.Range("E10").Value = 12.4203574294197
$ws.Range("G10").Value = 63.33793689455374
$ws.Range("H10").Value = 22.58312595171897
$ws.Range("I10").Value = 34.90264618215532
$ws.Range("J10").Value = 7.874345715386405
$ws.Range("L10").Value = 13.22277823072478
$ws.Range("M10").Value = 20.10465982297801
$ws.Range("B11").Value = 22.66955314599954
$ws.Range("C11").Value = 9.76389343348168
$ws.Range("D11").Value = 12.76605734223747
$ws.Range("E11").Value = 12.40627167333168
$ws.Range("G11").Value = 63.49605629075417
$ws.Range("H11").Value = 22.57552621452525
$ws.Range("I11").Value = 34.87995457694138
$ws.Range("J11").Value = 7.884570449475121
$ws.Range("L11").Value = 13.22061468917521
$ws.Range("M11").Value = 20.15593403723379
$ws.Range("B12").Value = 22.75881404418953
$ws.Range("C12").Value = 9.86367528553737
$ws.Range("D12").Value = 12.76486344733931
$ws.Range("E12").Value = 12.40104902762288
$ws.Range("G12").Value = 63.55838625879366
$ws.Range("H12").Value = 22.57330800603722
$ws.Range("I12").Value = 34.87242650684172
$ws.Range("J12").Value = 7.888409369377845
$ws.Range("L12").Value = 13.22007992693212
$ws.Range("M12").Value = 20.17584978425611
$ws.Range("B13").Value = 22.73958794896699
$ws.Range("C13").Value = 9.842260895165646
$ws.Range("D13").Value = 12.76510730397604
$ws.Range("E13").Value = 12.40216887370126
$ws.Range("G13").Value = 63.5448536464286
$ws.Range("H13").Value = 22.57375637559966
$ws.Range("I13").Value = 34.87400038917567
$ws.Range("J13").Value = 7.887584055703775
$ws.Range("L13").Value = 13.22018245379593
$ws.Range("M13").Value = 20.17153853829282
$ws.Range("B14").Value = 22.67689559902898
$ws.Range("C14").Value = 9.772136448100689
$ws.Range("D14").Value = 12.76595313964065
$ws.Range("E14").Value = 12.40583977455995
$ws.Range("G14").Value = 63.50113517883794
$ws.Range("H14").Value = 22.57533049240548
$ws.Range("I14").Value = 34.8793138746461
$ws.Range("J14").Value = 7.884886940495536
$ws.Range("L14").Value = 13.22056499544146
$ws.Range("M14").Value = 20.15756257302596
$ws.Range("B15").Value = 22.63850246338419
$ws.Range("C15").Value = 9.728963342748411
$ws.Range("D15").Value = 12.76651009944156
$ws.Range("E15").Value = 12.40810279188277
$ws.Range("G15").Value = 63.47467517649517
$ws.Range("H15").Value = 22.5763806341706
$ws.Range("I15").Value = 34.88270732601207
$ws.Range("J15").Value = 7.883230586010647
$ws.Range("L15").Value = 13.22083634729926
$ws.Range("M15").Value = 20.14906659422347
$ws.Range("B16").Value = 22.41866766684394
$ws.Range("C16").Value = 9.478230762920177
$ws.Range("D16").Value = 12.77029763381332
$ws.Range("E16").Value = 12.42129356915428
$ws.Range("G16").Value = 63.32794802952603
$ws.Range("H16").Value = 22.58371483681562
$ws.Range("I16").Value = 34.90427789271649
$ws.Range("J16").Value = 7.873672892947897
$ws.Range("L16").Value = 13.2229594730121
$ws.Range("M16").Value = 20.10137939019259
$ws.Range("B17").Value = 22.28403551701279
$ws.Range("C17").Value = 9.321499397300027
$ws.Range("D17").Value = 12.77315300394437
$ws.Range("E17").Value = 12.42958444281527
$ws.Range("G17").Value = 63.24233018053439
$ws.Range("H17").Value = 22.5893873676048
$ws.Range("I17").Value = 34.91940251644923
$ws.Range("J17").Value = 7.867750861451432
$ws.Range("L17").Value = 13.22476938583304
$ws.Range("M17").Value = 20.0730256062322
$ws.Range("B18").Value = 22.20671016707046
$ws.Range("C18").Value = 9.230302326887978
$ws.Range("D18").Value = 12.77499084513848
$ws.Range("E18").Value = 12.43442632100259
$ws.Range("G18").Value = 63.1947067625448
$ws.Range("H18").Value = 22.59308067062052
$ws.Range("I18").Value = 34.92879556260197
$ws.Range("J18").Value = 7.864323100083766
$ws.Range("L18").Value = 13.22599707157015
$ws.Range("M18").Value = 20.05705163146407
$ws.Range("B19").Value = 22.18055088067363
$ws.Range("C19").Value = 9.19924655625727
$ws.Range("D19").Value = 12.77564668963878
$ws.Range("E19").Value = 12.43607828135168
$ws.Range("G19").Value = 63.17886141495155
$ws.Range("H19").Value = 22.59440506268962
$ws.Range("I19").Value = 34.93209489810533
$ws.Range("J19").Value = 7.863158815568857
$ws.Range("L19").Value = 13.22644482265786
$ws.Range("M19").Value = 20.05170085648105
$ws.Range("B20").Value = 22.29835641587486
$ws.Range("C20").Value = 9.338292838488444
$ws.Range("D20").Value = 12.77282881174008
$ws.Range("E20").Value = 12.42869429384463
$ws.Range("G20").Value = 63.25127665108469
$ws.Range("H20").Value = 22.58873893747052
$ws.Range("I20").Value = 34.91772063921773
$ws.Range("J20").Value = 7.868383503862421
$ws.Range("L20").Value = 13.22455740179148
$ws.Range("M20").Value = 20.07600938290478
$ws.Range("B21").Value = 22.69530839705594
$ws.Range("C21").Value = 9.792779619354288
$ws.Range("D21").Value = 12.76569659908864
$ws.Range("E21").Value = 12.40475852421781
$ws.Range("G21").Value = 63.51390995439999
$ws.Range("H21").Value = 22.5748502217275
$ws.Range("I21").Value = 34.87772424438479
$ws.Range("J21").Value = 7.885680043414564
$ws.Range("L21").Value = 13.22044491702325
$ws.Range("M21").Value = 20.16165418824723
$ws.Range("B22").Value = 22.95515877104501
$ws.Range("C22").Value = 10.08002924825255
$ws.Range("D22").Value = 12.76277478943448
$ws.Range("E22").Value = 12.38976379486862
$ws.Range("G22").Value = 63.69984462662306
$ws.Range("H22").Value = 22.56961863654962
$ws.Range("I22").Value = 34.85779226099267
$ws.Range("J22").Value = 7.896792177624759
$ws.Range("L22").Value = 13.21941525055531
$ws.Range("M22").Value = 20.22053311408443
$ws.Range("B23").Value = 22.81646075324307
$ws.Range("C23").Value = 9.927633102968626
$ws.Range("D23").Value = 12.76417513996533
$ws.Range("E23").Value = 12.39770755624755
$ws.Range("G23").Value = 63.5993085191459
$ws.Range("H23").Value = 22.57205849608692
$ws.Range("I23").Value = 34.86786092504083
$ws.Range("J23").Value = 7.890879010603285
$ws.Range("L23").Value = 13.21981331041486
$ws.Range("M23").Value = 20.18884613179702
$ws.Range("B24").Value = 22.29188169285505
$ws.Range("C24").Value = 9.330703914089629
$ws.Range("D24").Value = 12.77297476773642
$ws.Range("E24").Value = 12.42909649547558
$ws.Range("G24").Value = 63.24722696808701
$ws.Range("H24").Value = 22.58903074690684
$ws.Range("I24").Value = 34.91847884295084
$ws.Range("J24").Value = 7.868097558204569
$ws.Range("L24").Value = 13.22465265680781
$ws.Range("M24").Value = 20.07465939812948
$ws.Range("B25").Value = 21.73065377138619
$ws.Range("C25").Value = 8.647077753138593
$ws.Range("D25").Value = 12.78994106596918
$ws.Range("E25").Value = 12.46575455068026
$ws.Range("G25").Value = 62.92950040733726
$ws.Range("H25").Value = 22.6237918085262
$ws.Range("I25").Value = 34.99958054725432
$ws.Range("J25").Value = 7.842699554496314
$ws.Range("L25").Value = 13.23701099610057
$ws.Range("M25").Value = 19.96431034052842

Write-Output "Updated cells"